# Add team record (Wins / Losses / Ties) columns to the BAL_2012 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 54

# --- Header row (row 1): new columns AD, AE, AF -------------------------
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Re-use the same cell format as the rest of the header row (e.g. AB1
# "Salary") so the new header cells look consistent with the existing ones.
$ws.Range("AB1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)   # xlPasteFormats

# --- Data rows (2-54): same team record for every player row -----------
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 93   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 69   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}
